# Weekly Fruta/Hortalizas update: a new price record for
# Vega Monumental Concepción - Arándano (blue) is inserted as row 86,
# pushing the previously existing rows 86-91 down to 87-92.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 86 (shifts old rows 86..91 -> 87..92)
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the latest weekly record
$ws.Range("A86").Value = 11
$ws.Range("B86").Value = "Vega Monumental Concepción"
$ws.Range("C86").Value = "Bíobío"
$ws.Range("D86").Value = 44610
$ws.Range("E86").Value = 8
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100101
$ws.Range("H86").Value = "Berries"
$ws.Range("I86").Value = 100101001
$ws.Range("J86").Value = "Arándano (blue)"
$ws.Range("K86").Value = "Sin especificar"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 180
$ws.Range("N86").Value = 2800
$ws.Range("O86").Value = 3000
$ws.Range("P86").Value = 2911
$ws.Range("Q86").Value = "$/bandeja 2 kilos"
$ws.Range("R86").Value = "Provincia de Curicó"
$ws.Range("S86").Value = 1456
$ws.Range("T86").Value = 2
